$d = $word.ActiveDocument

# --- 1) Fix the duplicated "les modalités de la rupture" phrase -----------
# The paragraph currently reads (across two runs):
#   "... nous pourrons déterminer d’un commun accord les modalités de la
#   rupture " + "les modalités de la rupture de mon contrat de travail."
# i.e. "les modalités de la rupture" is repeated. The fix drops the
# duplicated tail from the first run, which ends up split into two runs:
#   "déterminer d’un commun " and "accord ".
$hit = $d.Content
$found = $hit.Find.Execute("déterminer d’un commun accord les modalités de la rupture ", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    Write-Host "ERROR: target sentence not found"
}
$runStart = $hit.Start

$leadLen = ("déterminer d’un commun ").Length
$accordLen = ("accord ").Length

# Toggling a character property on just the "accord " slice forces Word to
# split the original run at that boundary (lead / accord / dup-tail) while
# leaving the run's actual formatting unchanged.
$accordRange = $d.Range($runStart + $leadLen, $runStart + $leadLen + $accordLen)
$accordRange.Bold = 1
$accordRange.Bold = 0

# Remove the now-isolated duplicate tail ("les modalités de la rupture ").
$dupRange = $d.Range($runStart + $leadLen + $accordLen, $hit.End)
$dupRange.Delete()

# --- 2) Move the _GoBack bookmark ------------------------------------------
# It currently wraps the "« Prénom Nom du salarié »" / "« Signature »"
# paragraphs at the end of the letter; it should instead sit (empty) in the
# blank paragraph right after the "Veuillez agréer, ..." sign-off paragraph.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Veuillez agréer*") {
        $targetPara = $d.Paragraphs.Item($i + 1)
        break
    }
}

# Insert a throwaway character so the empty paragraph has a run to anchor
# the bookmark to, wrap the bookmark around it, then delete the character -
# leaving an empty bookmarkStart/bookmarkEnd pair in that paragraph.
$insertAt = $d.Range($targetPara.Range.Start, $targetPara.Range.Start)
$insertAt.InsertBefore("X")
$anchorRange = $d.Range($targetPara.Range.Start, $targetPara.Range.Start + 1)
$d.Bookmarks.Add("_GoBack", $anchorRange)
$d.Bookmarks.Item("_GoBack").Range.Text = ""
